# SMDSDocuments.xlsx update: "updated SMDSDocument list for import"
#
# Net effect (per the target diff):
#  - Insert 2 new rows before the old row 467 (pushing the "Invest*" template
#    rows, previously at 467-480, down to 469-482).
#  - Re-purpose rows 465-468 with new/updated "Letter" entries:
#      465: Letter Redefer of Case
#      466: MED Procedural Order   (brand-new shared string)
#      467: Initial Letter without MED - CHD
#      468: Initial Letter without MED - CHG
#  - The two brand-new rows at the bottom (481-482) duplicate what is already
#    at rows 479-480 (Reconsideration Memo / Tabled Matter - Supplemental Memo),
#    matching the template row's normal layout (style ids preserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new rows at 467 - this shifts the existing rows 467:480
#    down to 469:482, carrying their values/styles with them automatically.
# ---------------------------------------------------------------------
$ws.Rows("467:468").Insert()

# Newly inserted rows lose their explicit row height; restore it to match
# the rest of the sheet (ht="15" customHeight="1").
$ws.Rows("467:468").RowHeight = 15

# ---------------------------------------------------------------------
# 2. Row 465: "Blurb" -> "Letter Redefer of Case"
# ---------------------------------------------------------------------
$ws.Range("C465").Value = "Letter Redefer of Case"
$ws.Range("D465").Value = "Letter Redefer of Case.docx"
$ws.Range("I465").Value = "NULL"

# ---------------------------------------------------------------------
# 3. Row 466: "Letter Redefer of Case" -> "MED Procedural Order" (new string)
# ---------------------------------------------------------------------
$ws.Range("C466").Value = "MED Procedural Order"
$ws.Range("D466").Value = "Med Procedural Order - All Types.docx"

# ---------------------------------------------------------------------
# 4. Row 467 (new, blank after insert): Initial Letter without MED - CHD
# ---------------------------------------------------------------------
$ws.Range("A467").Value = "ULP"
$ws.Range("B467").Value = "Letter"
$ws.Range("C467").Value = "Initial Letter without MED - CHD"
$ws.Range("D467").Value = "Initial Letter without MED - CHD.docx"
$ws.Range("E467").Value = 1
$ws.Range("F467").Value = "NULL"
$ws.Range("G467").Value = "NULL"
$ws.Range("H467").Value = "NULL"
$ws.Range("I467").Value = ""
$ws.Range("J467").Value = "CHD"
$ws.Range("K467").Value = "IL-CHD-Body.txt"
$ws.Range("L467").Value = "NULL"
$ws.Range("M467").Value = "NULL"
$ws.Range("N467").Value = "NULL"
$ws.Range("O467").Value = "NULL"

# ---------------------------------------------------------------------
# 5. Row 468 (new, blank after insert): Initial Letter without MED - CHG
# ---------------------------------------------------------------------
$ws.Range("A468").Value = "ULP"
$ws.Range("B468").Value = "Letter"
$ws.Range("C468").Value = "Initial Letter without MED - CHG"
$ws.Range("D468").Value = "Initial Letter without MED - CHG.docx"
$ws.Range("E468").Value = 1
$ws.Range("F468").Value = "NULL"
$ws.Range("G468").Value = "NULL"
$ws.Range("H468").Value = "NULL"
$ws.Range("I468").Value = ""
$ws.Range("J468").Value = "CHG"
$ws.Range("K468").Value = "IL-CHD-Body.txt"
$ws.Range("L468").Value = "NULL"
$ws.Range("M468").Value = "NULL"
$ws.Range("N468").Value = "NULL"
$ws.Range("O468").Value = "NULL"

# ---------------------------------------------------------------------
# 6. Column width tweaks (best effort - Excel quantizes widths to whole
#    pixels, so we target the values that round-trip closest to the
#    widths recorded in the target workbook).
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 8.1666666666667
$ws.Columns("G").ColumnWidth = 12.6666666666667
$ws.Columns("I").ColumnWidth = 20.1666666666667
$ws.Columns("J").ColumnWidth = 8.1666666666667
$ws.Columns("K").ColumnWidth = 21.8333333333333

# ---------------------------------------------------------------------
# 7. Update the view state (selection / scroll position) to match.
# ---------------------------------------------------------------------
$ws.Range("D448").Select()
$excel.ActiveWindow.ScrollRow = 448
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("M470").Select()
